# Automatische test-sync: 2025-08-01 23:54:50
# Appends the newest test-mail log entry (row 14) to the "Logs" sheet,
# extends the conditional formatting ranges to cover the new row, and
# bumps the matching "Overig" tally on the "Dashboard" sheet from 6 to 7.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 14

$logs.Cells.Item($newRow, 1).Value = "Hoi, hebben jullie al iets gehoord?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #9: Hoi, hebben jullie al iets gehoord?"
$logs.Cells.Item($newRow, 4).Value = "Overig"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-01 23:54:01"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Stretch every conditional-formatting block (columns D, G, H, I, J) one
# row further down so it keeps covering the whole data range through row 14.
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range($col + "2:" + $col + "13")
    $newRange = $logs.Range($col + "2:" + $col + "14")
    $fc = $oldRange.FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($newRange)
    }
}

# The new row belongs to the "Overig" category, so the Dashboard tally goes up by one.
$dashboard.Range("B2").Value = 7
